$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay text (e.g. "55.370.36"),
# so force text format before assigning to avoid Excel auto-numeric coercion.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.370.36'
$ws.Range("E2").Value = '  -4.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.945.66'
$ws.Range("E3").Value = '  -6.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.44'
$ws.Range("E5").Value = '  -7.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.90'
$ws.Range("E6").Value = '  -2.48%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.941.18'
$ws.Range("E8").Value = '  -6.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.416'
$ws.Range("E9").Value = '  -7.91%  '
$ws.Range("E10").Value = '  -4.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0995'
$ws.Range("E11").Value = '  -10.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.346'
$ws.Range("E12").Value = '  -11.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.448.61'
$ws.Range("E14").Value = '  -6.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.17'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '55.208.98'
$ws.Range("E16").Value = '  -4.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.942.46'
$ws.Range("E17").Value = '  -6.86%  '
$ws.Range("E18").Value = '  -9.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.56'
$ws.Range("E19").Value = '  -4.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.92'
$ws.Range("E20").Value = '  -8.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.37'
$ws.Range("E21").Value = '  -8.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '310.00'
$ws.Range("E22").Value = '  -10.60%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.458'
$ws.Range("E24").Value = '  -10.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '59.66'
$ws.Range("E25").Value = '  -14.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -5.66%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0842'
$ws.Range("E29").Value = '  -12.55%  '
$ws.Range("E30").Value = '  -5.49%  '
$ws.Range("E31").Value = '  -3.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.42'
$ws.Range("E32").Value = '  -7.07%  '
$ws.Range("E33").Value = '  -11.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.19'
$ws.Range("E34").Value = '  -11.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '146.94'
$ws.Range("E35").Value = '  -7.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.32'
$ws.Range("E36").Value = '  -12.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.60'
$ws.Range("E37").Value = '  -10.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.26'
$ws.Range("E38").Value = '  -10.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.53'
$ws.Range("E39").Value = '  -9.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0642'
$ws.Range("E40").Value = '  -7.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.974.19'
$ws.Range("E41").Value = '  -6.61%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.93'
$ws.Range("E43").Value = '  -11.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.991'
$ws.Range("E44").Value = '  -8.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.629'
$ws.Range("E45").Value = '  -10.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.36'
$ws.Range("E46").Value = '  -7.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.50'
$ws.Range("E47").Value = '  -11.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.111.00'
$ws.Range("E48").Value = '  -7.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0228'
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.77'
$ws.Range("E50").Value = '  -8.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.51'
$ws.Range("E51").Value = '  -10.98%  '
